$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit after "Vigenere
#    Cipher" (it gets relocated to the end of the new work-log entry below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fill in the blank work-log row (3rd data row of the single table) with
#    a date and a description of the day's work.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)

# --- Date cell -------------------------------------------------------------
$dateCell = $table.Cell(3, 1)
$dateRange = $dateCell.Range
$dateRange.Text = "28-04-15"

# --- Work cell ---------------------------------------------------------------
$workCell = $table.Cell(3, 2)
$workRange = $workCell.Range

# Trailing placeholder character - lets us park a collapsed Range exactly at
# the end of the paragraph text (the engine refuses to place a bookmark at
# the very last position of a story) and then get removed again afterwards.
$sentence1 = "Am started writing functionality of Vigenere_Algorithm."
$sentence2 = " First am getting message and key from user. Then am passing it to the "
$word3 = "vigenere"
$tail = " class and trying to encrypt the message using the key."
$fullText = $sentence1 + $sentence2 + $word3 + $tail + "#"

$workRange.Text = $fullText

$workCell2 = $table.Cell(3, 2)
$workRange2 = $workCell2.Range
$base = $workRange2.Start
$totalLen = $fullText.Length

# Apply the paragraph's font explicitly to the whole run of text, matching
# the Times New Roman formatting used throughout the rest of the table.
$wholeRange = $d.Range($base, $base + $totalLen)
$wholeRange.Font.Name = "Times New Roman"

# Force run breaks at the same boundaries Word's proofing marks would sit at
# (re-applying identical formatting to a sub-range splits the run there).
$b1 = "Am started writing functionality of ".Length
$b2 = $b1 + "Vigenere_".Length
$b3 = $b2 + "Algorithm".Length
$b4 = $b3 + ".".Length
$b5 = $b4 + $sentence2.Length
$b6 = $b5 + $word3.Length

$d.Range($base + $b1, $base + $b2).Font.Name = "Times New Roman"
$d.Range($base + $b2, $base + $b3).Font.Name = "Times New Roman"
$d.Range($base + $b3, $base + $b4).Font.Name = "Times New Roman"
$d.Range($base + $b5, $base + $b6).Font.Name = "Times New Roman"

# Place the "_GoBack" bookmark collapsed right before the trailing "#"
# placeholder, i.e. exactly after "...using the key." once the placeholder
# is removed.
$bmPos = $base + $totalLen - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# Remove the trailing placeholder character now that the bookmark anchors the
# correct position.
$workCell3 = $table.Cell(3, 2)
$workRange3 = $workCell3.Range
$placeholder = $d.Range($workRange3.End - 2, $workRange3.End - 1)
$placeholder.Delete()
